$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated currency quotations (Cotação) and recomputed dependent columns
# (Preço de Compra = Preço Original * Cotação, Preço de Venda = ROUND(Preço de Compra * Margem, 3))

$ws.Range("D2").Value = 5.5322
$ws.Range("E2").Value = 5532.144678
$ws.Range("G2").Value = 7745.003

$ws.Range("D3").Value = 6.336443575000001
$ws.Range("E3").Value = 28513.9960875
$ws.Range("G3").Value = 57027.992

$ws.Range("D4").Value = 5.5322
$ws.Range("E4").Value = 4978.924677999999
$ws.Range("G4").Value = 8464.172

$ws.Range("D5").Value = 5.5322
$ws.Range("E5").Value = 4420.2278
$ws.Range("G5").Value = 7514.387

$ws.Range("D6").Value = 6.336443575000001
$ws.Range("E6").Value = 19009.330725
$ws.Range("G6").Value = 36117.728

$ws.Range("D7").Value = 5.5322
$ws.Range("E7").Value = 2658.111456
$ws.Range("G7").Value = 5316.223

$ws.Range("D8").Value = 324.08
$ws.Range("E8").Value = 6481.599999999999
$ws.Range("G8").Value = 7453.84
